# TC-42.xlsx update — "Add files via upload"
#
# Re-labels the test case, replaces the "Pass:" test-data row with a
# generic "Contraseña" placeholder, clears the now-unused personal test
# data rows below it (Fecha/Dir/Ciudad/Estado/Cod postal/Country/Phone),
# and clears the leftover duplicate text in the final ("Step 7") row of
# the step-details table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Test Case ID: TC-01 -> TC-42
$ws.Range("B1").Value = "TC-42"

# Test Data table: "Pass: PruebaTC1!" -> "Contraseña"
$ws.Range("E11").Value = "Contraseña"

# Remaining personal test-data rows are no longer used — clear them out
$ws.Range("E12").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("E18").Value = ""

# Step 7 (row 32) of the Step Details / Expected Results table is cleared
$ws.Range("B32").Value = ""
$ws.Range("D32").Value = ""
